$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (Strike#) values replacing the old ones in column G, rows 2-17
$newValues = @(3, 1, 2, 2, 1, 0, 1, 1, 2, 0, 1, 2, 0, 1, 0, 1)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $newValues[$i]
}
